# Apply the changes described by the commit:
#   - Title font loses its 14pt size (drops back to the default size) and
#     becomes white; the header-row font also becomes white (bold stays).
#     (In the underlying styles.xml this merges what used to be two
#     separate bold fonts into a single bold+white font shared by the
#     title and the header row.)
#   - H3 (PERIOD TO EXPIRE) changes from -40 to -48
#   - I3 (LAST UPDATE) changes from 08-Sep-2025 to 16-Sep-2025 (kept as text)

$wb = $excel.ActiveWorkbook
$white = 16777215  # RGB(255,255,255)

foreach ($ws in $wb.Worksheets) {

    # --- Title cell (row 1, column A) -------------------------------------
    $title = $ws.Range("A1")
    $title.Font.Size = 11        # drop the 14pt size back to the default
    $title.Font.Color = $white   # make the title text white

    # --- Header row (row 2, spanning the used columns) ---------------------
    $lastCol = $ws.UsedRange.Columns.Count
    $headerRow = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRow.Font.Color = $white   # make header text white (stays bold)
}

# --- Training Dashboard data row updates -----------------------------------
$ws1 = $wb.Worksheets.Item("Training Dashboard")

# PERIOD TO EXPIRE: -40 -> -48
$ws1.Range("H3").Value = -48

# LAST UPDATE: 08-Sep-2025 -> 16-Sep-2025, stored as literal text (not a date)
$i3 = $ws1.Range("I3")
$i3.NumberFormat = "@"
$i3.Value = "16-Sep-2025"
